# Auto-generated edit script applying the Yojimbo_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53: No Accounting for Waste
$ws.Range("H53").Value = 245.77777
$ws.Range("I53").Value = 50
$ws.Range("J53").Value = 270.25
$ws.Range("K53").Value = 50
$ws.Range("L53").Value = 270.25
$ws.Range("M53").Value = 587
$ws.Range("N53").Value = -1544.25
# Row 58: A Matter of Vital Importance
$ws.Range("H58").Value = 553.26666
$ws.Range("I58").Value = 208.25
$ws.Range("J58").Value = 1933.3334
$ws.Range("K58").Value = 624.75
$ws.Range("L58").Value = 5800.0002
$ws.Range("M58").Value = -474.75
$ws.Range("N58").Value = -6100.0002
# Row 111: An Eye for Healing
$ws.Range("H111").Value = 353
$ws.Range("J111").Value = 333
$ws.Range("L111").Value = 999
$ws.Range("N111").Value = -7133
# Row 116: Growing Up
$ws.Range("H116").Value = 2953.3635
$ws.Range("I116").Value = 2544.0588
$ws.Range("J116").Value = 3388.25
$ws.Range("K116").Value = 2544.0588
$ws.Range("L116").Value = 3388.25
$ws.Range("M116").Value = 897.9412000000002
$ws.Range("N116").Value = -10272.25
# Row 129: Practical Command
$ws.Range("H129").Value = 559.4286
$ws.Range("J129").Value = 1083.5714
$ws.Range("L129").Value = 3250.7142
$ws.Range("N129").Value = -13250.7142
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1275.963
$ws.Range("I132").Value = 1341.1111
$ws.Range("J132").Value = 950.2222
$ws.Range("K132").Value = 4023.3333
$ws.Range("L132").Value = 2850.6666
$ws.Range("M132").Value = -1493.3333
$ws.Range("N132").Value = -7910.6666

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7873.59
$ws.Range("I32").Value = 2755.1743
$ws.Range("J32").Value = 39315.285
$ws.Range("K32").Value = 2755.1743
$ws.Range("L32").Value = 39315.285
$ws.Range("M32").Value = -2468.1743
$ws.Range("N32").Value = -39889.285
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 762.069
$ws.Range("I74").Value = 483.78946
$ws.Range("J74").Value = 1290.8
$ws.Range("K74").Value = 483.78946
$ws.Range("L74").Value = 1290.8
$ws.Range("M74").Value = 390.21054
$ws.Range("N74").Value = -3038.8
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 762.069
$ws.Range("I77").Value = 483.78946
$ws.Range("J77").Value = 1290.8
$ws.Range("K77").Value = 2418.9473
$ws.Range("L77").Value = 6454
$ws.Range("M77").Value = 1949.0527
$ws.Range("N77").Value = -15190
# Row 97: Ore for Me
$ws.Range("H97").Value = 748.4524
$ws.Range("I97").Value = 671.58826
$ws.Range("J97").Value = 1075.125
$ws.Range("K97").Value = 671.58826
$ws.Range("L97").Value = 1075.125
$ws.Range("M97").Value = -175.58826
$ws.Range("N97").Value = -2067.125
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2986.8823
$ws.Range("I122").Value = 2374.889
$ws.Range("J122").Value = 3675.375
$ws.Range("K122").Value = 7124.667
$ws.Range("L122").Value = 11026.125
$ws.Range("M122").Value = -4674.667
$ws.Range("N122").Value = -15926.125
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2364.827
$ws.Range("I132").Value = 1633.9286
$ws.Range("J132").Value = 5434.6
$ws.Range("K132").Value = 4901.7858
$ws.Range("L132").Value = 16303.8
$ws.Range("M132").Value = -2371.7858
$ws.Range("N132").Value = -21363.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 761.6667
$ws.Range("I16").Value = 660
$ws.Range("K16").Value = 660
$ws.Range("M16").Value = -373
# Row 31: Wall Not Found
$ws.Range("H31").Value = 52113.285
$ws.Range("I31").Value = 88572.164
$ws.Range("J31").Value = 3501.4443
$ws.Range("K31").Value = 88572.164
$ws.Range("L31").Value = 3501.4443
$ws.Range("M31").Value = -88277.164
$ws.Range("N31").Value = -4091.4443
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 52113.285
$ws.Range("I34").Value = 88572.164
$ws.Range("J34").Value = 3501.4443
$ws.Range("K34").Value = 88572.164
$ws.Range("L34").Value = 3501.4443
$ws.Range("M34").Value = -88370.164
$ws.Range("N34").Value = -3905.4443
# Row 113: Patient Patients
$ws.Range("H113").Value = 761.6667
$ws.Range("I113").Value = 660
$ws.Range("K113").Value = 660
$ws.Range("M113").Value = 1510

$ws = $wb.Worksheets.Item("CUL")
# Row 18: Fisher of Men
$ws.Range("H18").Value = 981.0909
$ws.Range("I18").Value = 836.5
$ws.Range("J18").Value = 1366.6666
$ws.Range("K18").Value = 2509.5
$ws.Range("L18").Value = 4099.9998
$ws.Range("M18").Value = -2340.5
$ws.Range("N18").Value = -4437.9998
# Row 56: Culture Club
$ws.Range("H56").Value = 4230.3335
$ws.Range("I56").Value = 4230.3335
$ws.Range("K56").Value = 4230.3335
$ws.Range("M56").Value = -3700.3335
# Row 107: Slippery Service
$ws.Range("H107").Value = 470.2
$ws.Range("J107").Value = 512.75
$ws.Range("L107").Value = 1538.25
$ws.Range("N107").Value = -5378.25

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1421.2667
$ws.Range("I97").Value = 1568.8889
$ws.Range("J97").Value = 1199.8334
$ws.Range("K97").Value = 1568.8889
$ws.Range("L97").Value = 1199.8334
$ws.Range("M97").Value = -1072.8889
$ws.Range("N97").Value = -2191.8334
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2425
$ws.Range("I113").Value = 3100
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 3100
$ws.Range("L113").Value = 400
$ws.Range("M113").Value = -930
$ws.Range("N113").Value = -4740
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 1670.7539
$ws.Range("I122").Value = 1328.434
$ws.Range("J122").Value = 3182.6667
$ws.Range("K122").Value = 3985.302
$ws.Range("L122").Value = 9548.000100000001
$ws.Range("M122").Value = -1535.302
$ws.Range("N122").Value = -14448.0001
# Row 132: On Board for Lar
$ws.Range("H132").Value = 1314.591
$ws.Range("I132").Value = 1150.2543
$ws.Range("J132").Value = 2699.7144
$ws.Range("K132").Value = 3450.7629
$ws.Range("L132").Value = 8099.1432
$ws.Range("M132").Value = -920.7629000000002
$ws.Range("N132").Value = -13159.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 1783.5294
$ws.Range("I7").Value = 1075
$ws.Range("J7").Value = 2001.5385
$ws.Range("K7").Value = 1075
$ws.Range("L7").Value = 2001.5385
$ws.Range("M7").Value = -963
$ws.Range("N7").Value = -2225.5385
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 1811.3077
$ws.Range("I61").Value = 1591.2307
$ws.Range("J61").Value = 2031.3846
$ws.Range("K61").Value = 1591.2307
$ws.Range("L61").Value = 2031.3846
$ws.Range("M61").Value = -1389.2307
$ws.Range("N61").Value = -2435.3846
# Row 113: Peace in Rest
$ws.Range("H113").Value = 1811.3077
$ws.Range("I113").Value = 1591.2307
$ws.Range("J113").Value = 2031.3846
$ws.Range("K113").Value = 1591.2307
$ws.Range("L113").Value = 2031.3846
$ws.Range("M113").Value = 578.7692999999999
$ws.Range("N113").Value = -6371.3846
# Row 126: Battered Books
$ws.Range("H126").Value = 1783.5294
$ws.Range("I126").Value = 1075
$ws.Range("J126").Value = 2001.5385
$ws.Range("K126").Value = 3225
$ws.Range("L126").Value = 6004.6155
$ws.Range("M126").Value = -755
$ws.Range("N126").Value = -10944.6155
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 2164.5469
$ws.Range("I132").Value = 1939.7858
$ws.Range("J132").Value = 3737.875
$ws.Range("K132").Value = 5819.357400000001
$ws.Range("L132").Value = 11213.625
$ws.Range("M132").Value = -3289.357400000001
$ws.Range("N132").Value = -16273.625
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 2021.3816
$ws.Range("I136").Value = 1404.2307
$ws.Range("J136").Value = 3358.5417
$ws.Range("K136").Value = 4212.6921
$ws.Range("L136").Value = 10075.6251
$ws.Range("M136").Value = -1662.6921
$ws.Range("N136").Value = -15175.6251

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 715143.4399999999
$ws.Range("I96").Value = 900
$ws.Range("J96").Value = 1667468
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 1667468
$ws.Range("M96").Value = 473
$ws.Range("N96").Value = -1670214
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 501910.8
$ws.Range("I122").Value = 626594.4399999999
$ws.Range("K122").Value = 1879783.32
$ws.Range("M122").Value = -1877333.32
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 508.12698
$ws.Range("I132").Value = 272.02274
$ws.Range("J132").Value = 1054.8948
$ws.Range("K132").Value = 816.06822
$ws.Range("L132").Value = 3164.6844
$ws.Range("M132").Value = 1713.93178
$ws.Range("N132").Value = -8224.6844
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 326.69864
$ws.Range("I136").Value = 250.13559
$ws.Range("J136").Value = 649.3570999999999
$ws.Range("K136").Value = 750.4067700000001
$ws.Range("L136").Value = 1948.0713
$ws.Range("M136").Value = 1799.59323
$ws.Range("N136").Value = -7048.0713
